$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename a few of the building category labels in column A.
# Row 1 stays "building_category" (header, unchanged)
$ws.Range("A2").Value = "Apartment block"   # was "Apartment"
$ws.Range("A3").Value = "House"             # was "SmallHouse"
$ws.Range("A8").Value = "Retail"            # was "Shop"
$ws.Range("A14").Value = "Storage repairs"  # was "StorageRepairs"

# Clear the stored cell selection/active cell so the saved view doesn't
# pin the selection to A2 (matches removal of <selection> in sheetView).
$ws.Range("A1").Select() | Out-Null
